$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (Changed) date column C for rows 2-6 from 45185 to 45204
$ws.Range("C2:C6").Value = 45204
